$wb = $excel.ActiveWorkbook

# --- Overview sheet: update the cached status text from
#     "Handoff transform failed" to "Ready for handoff" everywhere it is used.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"

# --- zh-cn sheet: file is now handed off, so fill in the handoff file link,
#     the handoff timestamp, and flip the handoff reason to "Include".
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B2").Value = "Ready for handoff"
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTest/oltest/blob/561fbe8993a614bd90a7ba39acf27ce7cc3d7768/e2e/7e073687-e8e9-4b70-aa2e-2cd31a43bdc0.d1ddc4969f60068365f8f84009c2e1a40cd74277.zh-cn.xlf", "", "", "7e073687-e8e9-4b70-aa2e-2cd31a43bdc0.d1ddc4969f60068365f8f84009c2e1a40cd74277.zh-cn.xlf")
$zh.Range("D2").Value = "2016-02-16 10:20:42"
$zh.Range("H2").Value = "Include"

# --- de-de sheet: same treatment, with the de-de handoff artifact.
$de = $wb.Worksheets.Item("de-de")
$de.Range("B2").Value = "Ready for handoff"
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTest/oltest/blob/561fbe8993a614bd90a7ba39acf27ce7cc3d7768/e2e/7e073687-e8e9-4b70-aa2e-2cd31a43bdc0.d1ddc4969f60068365f8f84009c2e1a40cd74277.de-de.xlf", "", "", "7e073687-e8e9-4b70-aa2e-2cd31a43bdc0.d1ddc4969f60068365f8f84009c2e1a40cd74277.de-de.xlf")
$de.Range("D2").Value = "2016-02-16 10:20:54"
$de.Range("H2").Value = "Include"
